$wb = $excel.ActiveWorkbook
$wsMonthly = $wb.Worksheets.Item("Monthly Data")
$wsAnnual  = $wb.Worksheets.Item("Annual Data")

# --- Update title / release-date / next-update strings (shared workbook strings) ---
# Both worksheets carry their own copies of these header cells, so update both.
$wsMonthly.Range("A2").Value = "June 2018 Monthly Energy Review"
$wsMonthly.Range("A6").Value = "Release Date: June 26, 2018"
$wsMonthly.Range("A7").Value = "Next Update: July 26, 2018"

$wsAnnual.Range("A2").Value = "June 2018 Monthly Energy Review"
$wsAnnual.Range("A6").Value = "Release Date: June 26, 2018"
$wsAnnual.Range("A7").Value = "Next Update: July 26, 2018"

# --- Revise existing monthly rows 529-542 (Monthly Data sheet) with updated EIA figures ---
$wsMonthly.Range("C529").Value = 75.262
$wsMonthly.Range("D529").Value = 74.912
$wsMonthly.Range("E529").Value = 150.174
$wsMonthly.Range("G529").Value = 1397.261
$wsMonthly.Range("H529").Value = 1651.675
$wsMonthly.Range("I529").Value = 3048.936
$wsMonthly.Range("J529").Value = 4377.419
$wsMonthly.Range("L529").Value = 62134.631
$wsMonthly.Range("M529").Value = 66662.224

$wsMonthly.Range("C530").Value = 75.387
$wsMonthly.Range("D530").Value = 75.036
$wsMonthly.Range("E530").Value = 150.423
$wsMonthly.Range("G530").Value = 1282.036
$wsMonthly.Range("H530").Value = 1755.363
$wsMonthly.Range("I530").Value = 3037.399
$wsMonthly.Range("J530").Value = 4398.844
$wsMonthly.Range("L530").Value = 50661.45
$wsMonthly.Range("M530").Value = 55210.717

$wsMonthly.Range("C531").Value = 74.003
$wsMonthly.Range("D531").Value = 73.658
$wsMonthly.Range("E531").Value = 147.661
$wsMonthly.Range("G531").Value = 1274.833
$wsMonthly.Range("H531").Value = 1770.31
$wsMonthly.Range("I531").Value = 3045.143
$wsMonthly.Range("J531").Value = 4478.8
$wsMonthly.Range("L531").Value = 39948.145
$wsMonthly.Range("M531").Value = 44574.606

$wsMonthly.Range("C532").Value = 45.553
$wsMonthly.Range("D532").Value = 28.657
$wsMonthly.Range("E532").Value = 74.21
$wsMonthly.Range("G532").Value = 1075.583
$wsMonthly.Range("H532").Value = 1750.917
$wsMonthly.Range("I532").Value = 2826.5
$wsMonthly.Range("J532").Value = 4150.531
$wsMonthly.Range("L532").Value = 39158.963
$wsMonthly.Range("M532").Value = 43383.704

$wsMonthly.Range("C533").Value = 36.542
$wsMonthly.Range("D533").Value = 22.989
$wsMonthly.Range("E533").Value = 59.531
$wsMonthly.Range("G533").Value = 1178.029
$wsMonthly.Range("H533").Value = 1656.568
$wsMonthly.Range("I533").Value = 2834.597
$wsMonthly.Range("J533").Value = 4201.467
$wsMonthly.Range("L533").Value = 45081.935
$wsMonthly.Range("M533").Value = 49342.933

$wsMonthly.Range("C534").Value = 46.166
$wsMonthly.Range("D534").Value = 29.043
$wsMonthly.Range("E534").Value = 75.209
$wsMonthly.Range("G534").Value = 1242.702
$wsMonthly.Range("H534").Value = 1578.086
$wsMonthly.Range("I534").Value = 2820.788
$wsMonthly.Range("J534").Value = 4225.606
$wsMonthly.Range("L534").Value = 63250.414
$wsMonthly.Range("M534").Value = 67551.229

$wsMonthly.Range("C535").Value = 46.143
$wsMonthly.Range("D535").Value = 17.383
$wsMonthly.Range("E535").Value = 63.526
$wsMonthly.Range("G535").Value = 1321.161
$wsMonthly.Range("H535").Value = 1514.584
$wsMonthly.Range("I535").Value = 2835.745
$wsMonthly.Range("J535").Value = 4268.285
$wsMonthly.Range("L535").Value = 74236.728
$wsMonthly.Range("M535").Value = 78568.539

$wsMonthly.Range("C536").Value = 49.413
$wsMonthly.Range("D536").Value = 18.615
$wsMonthly.Range("E536").Value = 68.028
$wsMonthly.Range("G536").Value = 1292.162
$wsMonthly.Range("H536").Value = 1529.738
$wsMonthly.Range("I536").Value = 2821.9
$wsMonthly.Range("J536").Value = 4216.578
$wsMonthly.Range("L536").Value = 73889.93
$wsMonthly.Range("M536").Value = 78174.536

$wsMonthly.Range("C537").Value = 49.607
$wsMonthly.Range("D537").Value = 18.688
$wsMonthly.Range("E537").Value = 68.295
$wsMonthly.Range("G537").Value = 1157.257
$wsMonthly.Range("H537").Value = 1668.331
$wsMonthly.Range("I537").Value = 2825.588
$wsMonthly.Range("J537").Value = 4161.387
$wsMonthly.Range("L537").Value = 62385.216
$wsMonthly.Range("M537").Value = 66614.898

$wsMonthly.Range("C538").Value = 49.856
$wsMonthly.Range("D538").Value = 37.991
$wsMonthly.Range("E538").Value = 87.847
$wsMonthly.Range("G538").Value = 1126.246
$wsMonthly.Range("H538").Value = 1782.495
$wsMonthly.Range("I538").Value = 2908.741
$wsMonthly.Range("J538").Value = 4243.411
$wsMonthly.Range("L538").Value = 54621.445
$wsMonthly.Range("M538").Value = 58952.703

$wsMonthly.Range("C539").Value = 59.538
$wsMonthly.Range("D539").Value = 45.368
$wsMonthly.Range("E539").Value = 104.906
$wsMonthly.Range("G539").Value = 1092.706
$wsMonthly.Range("H539").Value = 1830.459
$wsMonthly.Range("I539").Value = 2923.165
$wsMonthly.Range("J539").Value = 4249.133
$wsMonthly.Range("L539").Value = 48179.203
$wsMonthly.Range("M539").Value = 52533.242

$wsMonthly.Range("C540").Value = 75.425
$wsMonthly.Range("D540").Value = 57.474
$wsMonthly.Range("E540").Value = 132.899
$wsMonthly.Range("G540").Value = 1280.244
$wsMonthly.Range("H540").Value = 1640.041
$wsMonthly.Range("I540").Value = 2920.285
$wsMonthly.Range("J540").Value = 4362.034
$wsMonthly.Range("L540").Value = 65006.425
$wsMonthly.Range("M540").Value = 69501.358

$wsMonthly.Range("C541").Value = 66.184
$wsMonthly.Range("D541").Value = 72.161
$wsMonthly.Range("E541").Value = 138.345
$wsMonthly.Range("F541").Value = 1430.645
$wsMonthly.Range("G541").Value = 1290.391
$wsMonthly.Range("H541").Value = 1553.562
$wsMonthly.Range("I541").Value = 2843.953
$wsMonthly.Range("J541").Value = 4274.598
$wsMonthly.Range("L541").Value = 63547.714
$wsMonthly.Range("M541").Value = 67960.657

$wsMonthly.Range("C542").Value = 53.674
$wsMonthly.Range("D542").Value = 58.521
$wsMonthly.Range("E542").Value = 112.195
$wsMonthly.Range("F542").Value = 1367.727
$wsMonthly.Range("G542").Value = 1087.427
$wsMonthly.Range("H542").Value = 1766.993
$wsMonthly.Range("I542").Value = 2854.42
$wsMonthly.Range("J542").Value = 4222.147
$wsMonthly.Range("L542").Value = 47964.848
$wsMonthly.Range("M542").Value = 52299.19

# --- Add new monthly rows 543-555 (March 2017 - March 2018) to Monthly Data sheet ---
$wsMonthly.Range("A543").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A543").HorizontalAlignment = -4131
$wsMonthly.Range("A543").Value = 42795
$wsMonthly.Range("B543").Value = "Not Available"
$wsMonthly.Range("C543").Value = 58.423
$wsMonthly.Range("D543").Value = 63.699
$wsMonthly.Range("E543").Value = 122.122
$wsMonthly.Range("F543").Value = 1437.669
$wsMonthly.Range("G543").Value = 1172.172
$wsMonthly.Range("H543").Value = 1664.293
$wsMonthly.Range("I543").Value = 2836.465
$wsMonthly.Range("J543").Value = 4274.134
$wsMonthly.Range("K543").Value = 0
$wsMonthly.Range("L543").Value = 48825.958
$wsMonthly.Range("M543").Value = 53222.214

$wsMonthly.Range("A544").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A544").HorizontalAlignment = -4131
$wsMonthly.Range("A544").Value = 42826
$wsMonthly.Range("B544").Value = "Not Available"
$wsMonthly.Range("C544").Value = 40.473
$wsMonthly.Range("D544").Value = 24.9
$wsMonthly.Range("E544").Value = 65.373
$wsMonthly.Range("F544").Value = 1440.81
$wsMonthly.Range("G544").Value = 1067.605
$wsMonthly.Range("H544").Value = 1629.855
$wsMonthly.Range("I544").Value = 2697.46
$wsMonthly.Range("J544").Value = 4138.27
$wsMonthly.Range("K544").Value = 0
$wsMonthly.Range("L544").Value = 44323.847
$wsMonthly.Range("M544").Value = 48527.49

$wsMonthly.Range("A545").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A545").HorizontalAlignment = -4131
$wsMonthly.Range("A545").Value = 42856
$wsMonthly.Range("B545").Value = "Not Available"
$wsMonthly.Range("C545").Value = 39.962
$wsMonthly.Range("D545").Value = 24.586
$wsMonthly.Range("E545").Value = 64.548
$wsMonthly.Range("F545").Value = 1482.486
$wsMonthly.Range("G545").Value = 1098.283
$wsMonthly.Range("H545").Value = 1604.724
$wsMonthly.Range("I545").Value = 2703.007
$wsMonthly.Range("J545").Value = 4185.493
$wsMonthly.Range("K545").Value = 0
$wsMonthly.Range("L545").Value = 50926.005
$wsMonthly.Range("M545").Value = 55176.046

$wsMonthly.Range("A546").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A546").HorizontalAlignment = -4131
$wsMonthly.Range("A546").Value = 42887
$wsMonthly.Range("B546").Value = "Not Available"
$wsMonthly.Range("C546").Value = 45.507
$wsMonthly.Range("D546").Value = 27.997
$wsMonthly.Range("E546").Value = 73.504
$wsMonthly.Range("F546").Value = 1401.664
$wsMonthly.Range("G546").Value = 1094.108
$wsMonthly.Range("H546").Value = 1617.146
$wsMonthly.Range("I546").Value = 2711.254
$wsMonthly.Range("J546").Value = 4112.918
$wsMonthly.Range("K546").Value = 0
$wsMonthly.Range("L546").Value = 58951.924
$wsMonthly.Range("M546").Value = 63138.346

$wsMonthly.Range("A547").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A547").HorizontalAlignment = -4131
$wsMonthly.Range("A547").Value = 42917
$wsMonthly.Range("B547").Value = "Not Available"
$wsMonthly.Range("C547").Value = 53.309
$wsMonthly.Range("D547").Value = 16.546
$wsMonthly.Range("E547").Value = 69.855
$wsMonthly.Range("F547").Value = 1494.46
$wsMonthly.Range("G547").Value = 1047.123
$wsMonthly.Range("H547").Value = 1838.301
$wsMonthly.Range("I547").Value = 2885.424
$wsMonthly.Range("J547").Value = 4379.884
$wsMonthly.Range("K547").Value = 0
$wsMonthly.Range("L547").Value = 69900.111
$wsMonthly.Range("M547").Value = 74349.85

$wsMonthly.Range("A548").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A548").HorizontalAlignment = -4131
$wsMonthly.Range("A548").Value = 42948
$wsMonthly.Range("B548").Value = "Not Available"
$wsMonthly.Range("C548").Value = 48.549
$wsMonthly.Range("D548").Value = 15.069
$wsMonthly.Range("E548").Value = 63.618
$wsMonthly.Range("F548").Value = 1528.056
$wsMonthly.Range("G548").Value = 1064.994
$wsMonthly.Range("H548").Value = 1807.254
$wsMonthly.Range("I548").Value = 2872.248
$wsMonthly.Range("J548").Value = 4400.304
$wsMonthly.Range("K548").Value = 0
$wsMonthly.Range("L548").Value = 65933.994
$wsMonthly.Range("M548").Value = 70397.916

$wsMonthly.Range("A549").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A549").HorizontalAlignment = -4131
$wsMonthly.Range("A549").Value = 42979
$wsMonthly.Range("B549").Value = "Not Available"
$wsMonthly.Range("C549").Value = 47.069
$wsMonthly.Range("D549").Value = 14.609
$wsMonthly.Range("E549").Value = 61.678
$wsMonthly.Range("F549").Value = 1468.767
$wsMonthly.Range("G549").Value = 1030.015
$wsMonthly.Range("H549").Value = 1809.249
$wsMonthly.Range("I549").Value = 2839.264
$wsMonthly.Range("J549").Value = 4308.031
$wsMonthly.Range("K549").Value = 0
$wsMonthly.Range("L549").Value = 54779.784
$wsMonthly.Range("M549").Value = 59149.493

$wsMonthly.Range("A550").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A550").HorizontalAlignment = -4131
$wsMonthly.Range("A550").Value = 43009
$wsMonthly.Range("B550").Value = "Not Available"
$wsMonthly.Range("C550").Value = 42.669
$wsMonthly.Range("D550").Value = 37.703
$wsMonthly.Range("E550").Value = 80.372
$wsMonthly.Range("F550").Value = 1469.57
$wsMonthly.Range("G550").Value = 1149.103
$wsMonthly.Range("H550").Value = 1641.732
$wsMonthly.Range("I550").Value = 2790.835
$wsMonthly.Range("J550").Value = 4260.405
$wsMonthly.Range("K550").Value = 0
$wsMonthly.Range("L550").Value = 50214.467
$wsMonthly.Range("M550").Value = 54555.244

$wsMonthly.Range("A551").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A551").HorizontalAlignment = -4131
$wsMonthly.Range("A551").Value = 43040
$wsMonthly.Range("B551").Value = "Not Available"
$wsMonthly.Range("C551").Value = 49.578
$wsMonthly.Range("D551").Value = 43.808
$wsMonthly.Range("E551").Value = 93.386
$wsMonthly.Range("F551").Value = 1456.863
$wsMonthly.Range("G551").Value = 1142.217
$wsMonthly.Range("H551").Value = 1650.271
$wsMonthly.Range("I551").Value = 2792.488
$wsMonthly.Range("J551").Value = 4249.351
$wsMonthly.Range("K551").Value = 0
$wsMonthly.Range("L551").Value = 50992.13
$wsMonthly.Range("M551").Value = 55334.867

$wsMonthly.Range("A552").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A552").HorizontalAlignment = -4131
$wsMonthly.Range("A552").Value = 43070
$wsMonthly.Range("B552").Value = "Not Available"
$wsMonthly.Range("C552").Value = 61.661
$wsMonthly.Range("D552").Value = 54.486
$wsMonthly.Range("E552").Value = 116.147
$wsMonthly.Range("F552").Value = 1558.946
$wsMonthly.Range("G552").Value = 1180.736
$wsMonthly.Range("H552").Value = 1605.369
$wsMonthly.Range("I552").Value = 2786.105
$wsMonthly.Range("J552").Value = 4345.051
$wsMonthly.Range("K552").Value = 0
$wsMonthly.Range("L552").Value = 58388.345
$wsMonthly.Range("M552").Value = 62849.543

$wsMonthly.Range("A553").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A553").HorizontalAlignment = -4131
$wsMonthly.Range("A553").Value = 43101
$wsMonthly.Range("B553").Value = "Not Available"
$wsMonthly.Range("C553").Value = 68.861
$wsMonthly.Range("D553").Value = 35.059
$wsMonthly.Range("E553").Value = 103.92
$wsMonthly.Range("F553").Value = 1689.106
$wsMonthly.Range("G553").Value = 1269.77
$wsMonthly.Range("H553").Value = 1667.075
$wsMonthly.Range("I553").Value = 2936.845
$wsMonthly.Range("J553").Value = 4625.951
$wsMonthly.Range("K553").Value = 0
$wsMonthly.Range("L553").Value = 64650.176
$wsMonthly.Range("M553").Value = 69380.047

$wsMonthly.Range("A554").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A554").HorizontalAlignment = -4131
$wsMonthly.Range("A554").Value = 43132
$wsMonthly.Range("B554").Value = "Not Available"
$wsMonthly.Range("C554").Value = 53.112
$wsMonthly.Range("D554").Value = 50.81
$wsMonthly.Range("E554").Value = 103.922
$wsMonthly.Range("F554").Value = 1388.187
$wsMonthly.Range("G554").Value = 1131.822
$wsMonthly.Range("H554").Value = 1729.148
$wsMonthly.Range("I554").Value = 2860.97
$wsMonthly.Range("J554").Value = 4249.157
$wsMonthly.Range("K554").Value = 0
$wsMonthly.Range("L554").Value = 45823.067
$wsMonthly.Range("M554").Value = 50176.146

$wsMonthly.Range("A555").NumberFormat = "yyyy mmmm"
$wsMonthly.Range("A555").HorizontalAlignment = -4131
$wsMonthly.Range("A555").Value = 43160
$wsMonthly.Range("B555").Value = "Not Available"
$wsMonthly.Range("C555").Value = 50.608
$wsMonthly.Range("D555").Value = 7.14
$wsMonthly.Range("E555").Value = 57.748
$wsMonthly.Range("F555").Value = 1113.576
$wsMonthly.Range("G555").Value = 1169.107
$wsMonthly.Range("H555").Value = 1610.478
$wsMonthly.Range("I555").Value = 2779.585
$wsMonthly.Range("J555").Value = 3893.161
$wsMonthly.Range("K555").Value = 0
$wsMonthly.Range("L555").Value = 44495.503
$wsMonthly.Range("M555").Value = 48446.411

# --- Revise Annual Data row 80 (2016 total) with updated EIA figures ---
$wsAnnual.Range("C80").Value = 682.895
$wsAnnual.Range("D80").Value = 499.814
$wsAnnual.Range("E80").Value = 1182.709
$wsAnnual.Range("G80").Value = 14720.22
$wsAnnual.Range("H80").Value = 20128.567
$wsAnnual.Range("I80").Value = 34848.787
$wsAnnual.Range("J80").Value = 51333.495
$wsAnnual.Range("L80").Value = 678554.486
$wsAnnual.Range("M80").Value = 731070.69

# --- Add new Annual Data row 81 (2017 total) ---
$wsAnnual.Range("A81").HorizontalAlignment = -4131
$wsAnnual.Range("A81").Value = 2017
$wsAnnual.Range("B81").Value = "Not Available"
$wsAnnual.Range("C81").Value = 607.058
$wsAnnual.Range("D81").Value = 454.085
$wsAnnual.Range("E81").Value = 1061.143
$wsAnnual.Range("F81").Value = 17537.663
$wsAnnual.Range("G81").Value = 13424.174
$wsAnnual.Range("H81").Value = 20188.749
$wsAnnual.Range("I81").Value = 33612.923
$wsAnnual.Range("J81").Value = 51150.586
$wsAnnual.Range("K81").Value = 0
$wsAnnual.Range("L81").Value = 664749.129
$wsAnnual.Range("M81").Value = 716960.858
